$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "MC2 Link" column (C2:C4) label text from "LINK" to "Click Me"
$ws.Range("C2").Value = "Click Me"
$ws.Range("C3").Value = "Click Me"
$ws.Range("C4").Value = "Click Me"

# Update the "Implemented by Capabilities" list in D3:
#  - rename COE-PE-* entries to COE-PD-*
#  - append new reference links to the end of the list
$ws.Range("D3").Value = "EB-UO-2B_Business_Product,EB-UO-LOTD_Operation,EB-UO-Subscription_Product_&_management,EB-ND-Network_Development,EB-Enabling-Market_Planning,EB-Enabling-Project_Management,EB-Enabling-Training,COE-HR-EHS,EB-SO-Service_Planning,EB-SO-Service_Operation,EB-SO-Service_Parts_&_Logistics_,EB-UO-Business_Intelligence,EB-Community-IP_Management,EB-BD-GR,COE-Finance-Tax,COE-Finance-Accounting,COE-PD-Product_Planning,COE-Purchase-Vendor_Management,COE-Purchase-Purchase_Operation,COE-HR-Culture_and_Internal_Communications,COE-HR-Organization_Development,COE-HR-Talent_Aquisition,COE-HR-Total_Rewards,COE-PD-Product_Experience,COE-NIO_Life-NIO_Life_Operation,COE-NIO_Life-NIO_Life_Product_Launch,COE-Power-Power_Operation,COE-Power-R_&_M,COE-Legal-EU_Legal,COE-Digital-Digital_Product_Planning,COE-Digital-Digital_Product_Design,COE-Digital-Digital_Product_Launch,COE-PD-Vehicle_Quality,EB-UO-Commercial_Product,EB-BD-Retail_Financing,EB-BD-LeaseCo_Development,COE-Digital-Digital_System_Operation,COE-Finance-Finance_Admin,COE-Finance-Finance_Planning,COE-Power-Battery_&_Safety,COE-Power-Market_Management,COE-Finance-Controling,EB-UD-PR,EB-UD-Communication,EB-UD-Events,EB-UO-Sales_Planning,EB-UO-Production_and_Order_Planning,EB-UO-Fleet_Planning,COE-HR-Learning_Development,COE-HR-HRBP_HQ,EB-HR-HRBP_HQ,EB-HR-HRBP_Country,COE-HR_Employee_Experience_and_Delivery,COE-Power-LTO,EB-UO-Partner_Strategy"
